$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.410.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3971"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08036"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.05%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.891.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.045"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.270"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001044"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06631"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.416.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.479"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.254"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.108.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.125"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.24%  "
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9753"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09505"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.590"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  +5.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.356"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06112"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02255"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.87%  "
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5958"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1883"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5596"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("E48").Value = "  +5.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06943"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.063"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.38%  "
